# Update scripts with new TPM values for the Lgi1-Rtn4r (YoungD7) sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MuSCs -> Lgi1 -> Rtn4r -> FAPs) gets refreshed with the new TPM
# derived numbers that used to live on row 3. The old ECs target-cluster row
# is dropped entirely, and the FAPs row takes its place as row 2.
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 0.1038603333333333
$ws.Range("H2").Value = 0.311581
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.6668756666666665
$ws.Range("N2").Value = 2.000627
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.06926192903188888
$ws.Range("R2").Value = 0.6233573612869999
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Row 3 (the ECs target-cluster row) is no longer present in the new export.
$ws.Rows(3).Delete()
